$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cell, $value)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

$ws.Cells.Item(2, 4).Value = "30.619.85"
$ws.Cells.Item(2, 5).Value = "  -1.57%  "
$ws.Cells.Item(3, 4).Value = "1.884.21"
$ws.Cells.Item(3, 5).Value = "  -1.65%  "
Set-TextValue $ws.Cells.Item(4, 4) "1.002"
$ws.Cells.Item(4, 5).Value = "  +0.25%  "
Set-TextValue $ws.Cells.Item(5, 4) "235.25"
$ws.Cells.Item(5, 5).Value = "  -4.29%  "
Set-TextValue $ws.Cells.Item(6, 4) "1.001"
$ws.Cells.Item(6, 5).Value = "  +0.11%  "
Set-TextValue $ws.Cells.Item(7, 4) "0.4878"
$ws.Cells.Item(7, 5).Value = "  -2.27%  "
Set-TextValue $ws.Cells.Item(8, 4) "0.2878"
$ws.Cells.Item(8, 5).Value = "  -4.15%  "
Set-TextValue $ws.Cells.Item(9, 4) "0.06653"
$ws.Cells.Item(9, 5).Value = "  -4.20%  "
$ws.Cells.Item(10, 4).Value = "1.880.15"
$ws.Cells.Item(10, 5).Value = "  -1.88%  "
Set-TextValue $ws.Cells.Item(11, 4) "16.77"
$ws.Cells.Item(11, 5).Value = "  -1.23%  "
Set-TextValue $ws.Cells.Item(12, 4) "0.07222"
$ws.Cells.Item(12, 5).Value = "  -1.35%  "
Set-TextValue $ws.Cells.Item(13, 4) "88.62"
$ws.Cells.Item(13, 5).Value = "  -1.15%  "
Set-TextValue $ws.Cells.Item(14, 4) "4.998"
Set-TextValue $ws.Cells.Item(15, 4) "0.6620"
$ws.Cells.Item(15, 5).Value = "  -3.08%  "
$ws.Cells.Item(16, 4).Value = "30.572.23"
Set-TextValue $ws.Cells.Item(17, 4) "0.000007820"
$ws.Cells.Item(17, 5).Value = "  -3.34%  "
$ws.Cells.Item(18, 5).Value = "  +0.09%  "
Set-TextValue $ws.Cells.Item(19, 4) "12.97"
$ws.Cells.Item(19, 5).Value = "  -3.49%  "
$ws.Cells.Item(20, 4).Value = "2.120.88"
$ws.Cells.Item(20, 5).Value = "  -1.77%  "
Set-TextValue $ws.Cells.Item(21, 4) "1.002"
$ws.Cells.Item(21, 5).Value = "  +0.29%  "
Set-TextValue $ws.Cells.Item(22, 4) "4.730"
$ws.Cells.Item(22, 5).Value = "  -3.18%  "
Set-TextValue $ws.Cells.Item(23, 4) "186.03"
$ws.Cells.Item(23, 5).Value = "  +6.07%  "
Set-TextValue $ws.Cells.Item(24, 4) "6.047"
$ws.Cells.Item(24, 5).Value = "  -0.55%  "
Set-TextValue $ws.Cells.Item(25, 4) "9.269"
Set-TextValue $ws.Cells.Item(26, 4) "157.29"
$ws.Cells.Item(26, 5).Value = "  +3.61%  "
Set-TextValue $ws.Cells.Item(27, 4) "18.25"
$ws.Cells.Item(27, 5).Value = "  +0.74%  "
Set-TextValue $ws.Cells.Item(28, 4) "1.832"
$ws.Cells.Item(28, 5).Value = "  -6.10%  "
Set-TextValue $ws.Cells.Item(29, 4) "1.402"
$ws.Cells.Item(29, 5).Value = "  -0.81%  "
Set-TextValue $ws.Cells.Item(30, 4) "4.251"
$ws.Cells.Item(30, 5).Value = "  -2.59%  "
Set-TextValue $ws.Cells.Item(31, 4) "0.09014"
$ws.Cells.Item(31, 5).Value = "  +0.77%  "
Set-TextValue $ws.Cells.Item(32, 4) "3.927"
$ws.Cells.Item(32, 5).Value = "  -3.37%  "
Set-TextValue $ws.Cells.Item(33, 4) "0.05195"
$ws.Cells.Item(33, 5).Value = "  -1.13%  "
Set-TextValue $ws.Cells.Item(34, 4) "0.7317"
$ws.Cells.Item(34, 5).Value = "  -2.27%  "
Set-TextValue $ws.Cells.Item(35, 4) "1.077"
$ws.Cells.Item(35, 5).Value = "  -5.92%  "
Set-TextValue $ws.Cells.Item(36, 4) "2.696"
$ws.Cells.Item(36, 5).Value = "  +1.22%  "
Set-TextValue $ws.Cells.Item(37, 4) "0.01813"
$ws.Cells.Item(37, 5).Value = "  -5.33%  "
Set-TextValue $ws.Cells.Item(38, 4) "2.654"
$ws.Cells.Item(38, 5).Value = "  -3.17%  "
Set-TextValue $ws.Cells.Item(39, 4) "0.9200"
$ws.Cells.Item(39, 5).Value = "  -2.41%  "
Set-TextValue $ws.Cells.Item(40, 4) "2.032"
$ws.Cells.Item(40, 5).Value = "  -7.66%  "
Set-TextValue $ws.Cells.Item(41, 4) "0.4302"
$ws.Cells.Item(41, 5).Value = "  -1.43%  "
Set-TextValue $ws.Cells.Item(42, 4) "104.23"
$ws.Cells.Item(42, 5).Value = "  -0.93%  "
Set-TextValue $ws.Cells.Item(43, 4) "0.9975"
$ws.Cells.Item(43, 5).Value = "  -0.25%  "
Set-TextValue $ws.Cells.Item(44, 4) "5.711"
$ws.Cells.Item(44, 5).Value = "  -4.09%  "
Set-TextValue $ws.Cells.Item(45, 4) "0.1340"
$ws.Cells.Item(45, 5).Value = "  +0.62%  "
Set-TextValue $ws.Cells.Item(46, 4) "7.271"
$ws.Cells.Item(46, 5).Value = "  -7.10%  "
Set-TextValue $ws.Cells.Item(47, 4) "0.05821"
$ws.Cells.Item(47, 5).Value = "  -0.63%  "
Set-TextValue $ws.Cells.Item(48, 4) "8.685"
$ws.Cells.Item(48, 5).Value = "  +0.65%  "
Set-TextValue $ws.Cells.Item(51, 4) "33.16"
$ws.Cells.Item(51, 5).Value = "  -0.61%  "

# Rows 49 and 50 swap content (Decentraland <-> NEARProtocol) with updated price/volume
$ws.Cells.Item(49, 2).Value = "Decentraland"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
Set-TextValue $ws.Cells.Item(49, 4) "0.3918"
$ws.Cells.Item(49, 5).Value = "  +0.55%  "
$ws.Cells.Item(50, 2).Value = "NEARProtocol"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue $ws.Cells.Item(50, 4) "1.413"
$ws.Cells.Item(50, 5).Value = "  +1.55%  "
